$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 37124.5
$ws.Range("J3").Value = 37124.5
$ws.Range("L3").Value = 37124.5
$ws.Range("N3").Value = -37352.5

$ws.Range("H87").Value = 35069.2
$ws.Range("J87").Value = 35069.2
$ws.Range("L87").Value = 35069.2
$ws.Range("N87").Value = -37565.2

$ws.Range("H90").Value = 35069.2
$ws.Range("J90").Value = 35069.2
$ws.Range("L90").Value = 105207.6
$ws.Range("N90").Value = -117687.6

$ws.Range("H93").Value = 39146.855
$ws.Range("J93").Value = 39146.855
$ws.Range("L93").Value = 39146.855
$ws.Range("N93").Value = -44138.855

$ws.Range("H95").Value = 37970.668
$ws.Range("J95").Value = 37970.668
$ws.Range("L95").Value = 37970.668
$ws.Range("N95").Value = -43462.668

$ws.Range("H102").Value = 37124.5
$ws.Range("J102").Value = 37124.5
$ws.Range("L102").Value = 37124.5
$ws.Range("N102").Value = -43614.5

$ws.Range("H105").Value = 47992
$ws.Range("J105").Value = 47992
$ws.Range("L105").Value = 47992
$ws.Range("N105").Value = -54980

$ws.Range("H112").Value = 1137964.9
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 1390623.6
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 4171870.8
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -4174086.8

$ws.Range("H127").Value = 1065.1666
$ws.Range("I127").Value = 488.33334
$ws.Range("J127").Value = 1642
$ws.Range("K127").Value = 1465.00002
$ws.Range("L127").Value = 4926
$ws.Range("M127").Value = 3494.99998
$ws.Range("N127").Value = -14846

$ws.Range("H132").Value = 24266.879
$ws.Range("I132").Value = 3293.8157
$ws.Range("K132").Value = 9881.447100000001
$ws.Range("M132").Value = -7351.447100000001

$ws.Range("H137").Value = 4228.3486
$ws.Range("J137").Value = 8192.723
$ws.Range("L137").Value = 24578.169
$ws.Range("N137").Value = -29678.169

$ws.Range("H138").Value = 1487.9166
$ws.Range("I138").Value = 1231.2693
$ws.Range("J138").Value = 3156.125
$ws.Range("K138").Value = 3693.8079
$ws.Range("L138").Value = 9468.375
$ws.Range("M138").Value = 1446.1921
$ws.Range("N138").Value = -19748.375

$ws.Range("H141").Value = 1522.3572
$ws.Range("I141").Value = 998.78845
$ws.Range("J141").Value = 8328.75
$ws.Range("K141").Value = 2996.36535
$ws.Range("L141").Value = 24986.25
$ws.Range("M141").Value = 2183.63465
$ws.Range("N141").Value = -35346.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11973.129
$ws.Range("I32").Value = 10648.519
$ws.Range("J32").Value = 20914.25
$ws.Range("K32").Value = 10648.519
$ws.Range("L32").Value = 20914.25
$ws.Range("M32").Value = -10361.519
$ws.Range("N32").Value = -21488.25

$ws.Range("H61").Value = 2387.3518
$ws.Range("I61").Value = 1489.3125
$ws.Range("K61").Value = 1489.3125
$ws.Range("M61").Value = -1277.3125

$ws.Range("H74").Value = 1861.9166
$ws.Range("I74").Value = 1249.6
$ws.Range("J74").Value = 2299.2856
$ws.Range("K74").Value = 1249.6
$ws.Range("L74").Value = 2299.2856
$ws.Range("M74").Value = -375.5999999999999
$ws.Range("N74").Value = -4047.2856

$ws.Range("H77").Value = 1861.9166
$ws.Range("I77").Value = 1249.6
$ws.Range("J77").Value = 2299.2856
$ws.Range("K77").Value = 6248
$ws.Range("L77").Value = 11496.428
$ws.Range("M77").Value = -1880
$ws.Range("N77").Value = -20232.428

$ws.Range("H97").Value = 766.2069
$ws.Range("I97").Value = 668.5238000000001
$ws.Range("J97").Value = 1022.625
$ws.Range("K97").Value = 668.5238000000001
$ws.Range("L97").Value = 1022.625
$ws.Range("M97").Value = -172.5238000000001
$ws.Range("N97").Value = -2014.625

$ws.Range("H110").Value = 1479.3667
$ws.Range("I110").Value = 1436.3704
$ws.Range("J110").Value = 1866.3334
$ws.Range("K110").Value = 1436.3704
$ws.Range("L110").Value = 1866.3334
$ws.Range("M110").Value = 608.6296
$ws.Range("N110").Value = -5956.3334

$ws.Range("H122").Value = 1935.697
$ws.Range("I122").Value = 2159.913
$ws.Range("K122").Value = 6479.739
$ws.Range("M122").Value = -4029.739

$ws.Range("H136").Value = 2387.3518
$ws.Range("I136").Value = 1489.3125
$ws.Range("K136").Value = 4467.9375
$ws.Range("M136").Value = -1917.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2194.7632
$ws.Range("I105").Value = 2348.75
$ws.Range("J105").Value = 2153.7
$ws.Range("K105").Value = 2348.75
$ws.Range("L105").Value = 2153.7
$ws.Range("M105").Value = -601.75
$ws.Range("N105").Value = -5647.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1900
$ws.Range("I16").Value = 1350
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1350
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1063
$ws.Range("N16").Value = -3574

$ws.Range("H31").Value = 182423.45
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 182423.45
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 182423.45
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -183013.45

$ws.Range("H34").Value = 182423.45
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 182423.45
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 182423.45
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -182827.45

$ws.Range("H105").Value = 3708.1667
$ws.Range("I105").Value = 3799.8
$ws.Range("J105").Value = 3250
$ws.Range("K105").Value = 3799.8
$ws.Range("L105").Value = 3250
$ws.Range("M105").Value = -2052.8
$ws.Range("N105").Value = -6744

$ws.Range("H113").Value = 1900
$ws.Range("I113").Value = 1350
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1350
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 820
$ws.Range("N113").Value = -7340

$ws.Range("H124").Value = 39996
$ws.Range("J124").Value = 39996
$ws.Range("L124").Value = 39996
$ws.Range("N124").Value = -44906

$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920

$ws.Range("H131").Value = 38318
$ws.Range("J131").Value = 38318
$ws.Range("L131").Value = 38318
$ws.Range("N131").Value = -48398

$ws.Range("H132").Value = 23820.889
$ws.Range("I132").Value = 1431.48
$ws.Range("J132").Value = 109934
$ws.Range("K132").Value = 4294.440000000001
$ws.Range("L132").Value = 329802
$ws.Range("M132").Value = -1764.440000000001
$ws.Range("N132").Value = -334862

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 30311346
$ws.Range("I137").Value = 2624.6155
$ws.Range("J137").Value = 50012016
$ws.Range("K137").Value = 7873.8465
$ws.Range("L137").Value = 150036048
$ws.Range("M137").Value = -2773.8465
$ws.Range("N137").Value = -150046248

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1683.4546
$ws.Range("I102").Value = 1721.8
$ws.Range("K102").Value = 1721.8
$ws.Range("M102").Value = -99.79999999999995

$ws.Range("H118").Value = 38298
$ws.Range("J118").Value = 38298
$ws.Range("L118").Value = 38298
$ws.Range("N118").Value = -41612

$ws.Range("H120").Value = 38317
$ws.Range("J120").Value = 38317
$ws.Range("L120").Value = 38317
$ws.Range("N120").Value = -47993

$ws.Range("H122").Value = 1543.3334
$ws.Range("I122").Value = 1698.5714
$ws.Range("K122").Value = 5095.7142
$ws.Range("M122").Value = -2645.7142

$ws.Range("H125").Value = 31481.334
$ws.Range("J125").Value = 31481.334
$ws.Range("L125").Value = 31481.334
$ws.Range("N125").Value = -36401.334

$ws.Range("H126").Value = 9586.571
$ws.Range("I126").Value = 15251.5
$ws.Range("J126").Value = 2033.3334
$ws.Range("K126").Value = 45754.5
$ws.Range("L126").Value = 6100.0002
$ws.Range("M126").Value = -43284.5
$ws.Range("N126").Value = -11040.0002

$ws.Range("H127").Value = 47303
$ws.Range("J127").Value = 47303
$ws.Range("L127").Value = 47303
$ws.Range("N127").Value = -57223

$ws.Range("H131").Value = 42318
$ws.Range("J131").Value = 42318
$ws.Range("L131").Value = 42318
$ws.Range("N131").Value = -52398

$ws.Range("H132").Value = 2529.9348
$ws.Range("I132").Value = 1629.96
$ws.Range("J132").Value = 3601.3333
$ws.Range("K132").Value = 4889.88
$ws.Range("L132").Value = 10803.9999
$ws.Range("M132").Value = -2359.88
$ws.Range("N132").Value = -15863.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 35281
$ws.Range("J109").Value = 35281
$ws.Range("L109").Value = 35281
$ws.Range("N109").Value = -38055

$ws.Range("H117").Value = 40134
$ws.Range("J117").Value = 40134
$ws.Range("L117").Value = 40134
$ws.Range("N117").Value = -49312

$ws.Range("H123").Value = 25130
$ws.Range("J123").Value = 29412.5
$ws.Range("L123").Value = 29412.5
$ws.Range("N123").Value = -39212.5

$ws.Range("H129").Value = 36616.668
$ws.Range("J129").Value = 36616.668
$ws.Range("L129").Value = 36616.668
$ws.Range("N129").Value = -46616.668

$ws.Range("H131").Value = 43326
$ws.Range("J131").Value = 43326
$ws.Range("L131").Value = 43326
$ws.Range("N131").Value = -53406

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 37993.6
$ws.Range("J115").Value = 37993.6
$ws.Range("L115").Value = 37993.6
$ws.Range("N115").Value = -41127.6

$ws.Range("H118").Value = 38996.445
$ws.Range("J118").Value = 41996
$ws.Range("L118").Value = 41996
$ws.Range("N118").Value = -45310

$ws.Range("H132").Value = 1619.8667
$ws.Range("I132").Value = 709.5
$ws.Range("J132").Value = 2985.4167
$ws.Range("K132").Value = 2128.5
$ws.Range("L132").Value = 8956.250100000001
$ws.Range("M132").Value = 401.5
$ws.Range("N132").Value = -14016.2501

$ws.Range("H136").Value = 18130.465
$ws.Range("I136").Value = 23822.697
$ws.Range("K136").Value = 71468.091
$ws.Range("M136").Value = -68918.091
